$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values row by row (each row shifted left by one column vs. previous values,
# reflecting the addition of the newest quarter column and removal of the oldest one).
$ws.Range("B2").Value = 0.1593188804880037
$ws.Range("C2").Value = -0.549633044125851
$ws.Range("D2").Value = -0.05752705180337387
$ws.Range("E2").Value = 0.4122776954696132
$ws.Range("F2").Value = 0.7251028904350592
$ws.Range("G2").Value = -0.1532996908165208
$ws.Range("H2").Value = -0.8832117059949898
$ws.Range("I2").Value = 0.7789673903946376
$ws.Range("J2").Value = 0.1597481019993938
$ws.Range("K2").Value = 0.3970518656191074
$ws.Range("B3").Value = -0.5671871654858154
$ws.Range("C3").Value = -0.07508117316333833
$ws.Range("D3").Value = 0.3947235741096488
$ws.Range("E3").Value = 0.7075487690750948
$ws.Range("F3").Value = -0.1708538121764852
$ws.Range("G3").Value = -0.9007658273549541
$ws.Range("H3").Value = 0.7614132690346732
$ws.Range("I3").Value = 0.1421939806394294
$ws.Range("J3").Value = 0.379497744259143
$ws.Range("K3").Value = 0.2762203743837313
$ws.Range("B4").Value = 0.01302614578766791
$ws.Range("C4").Value = 0.482830893060655
$ws.Range("D4").Value = 0.795656088026101
$ws.Range("E4").Value = -0.082746493225479
$ws.Range("F4").Value = -0.8126585084039479
$ws.Range("G4").Value = 0.8495205879856794
$ws.Range("H4").Value = 0.2303012995904356
$ws.Range("I4").Value = 0.4676050632101492
$ws.Range("J4").Value = 0.3643276933347375
$ws.Range("K4").Value = -0.4272707339946972
$ws.Range("B5").Value = 0.6614055265484386
$ws.Range("C5").Value = 0.9742307215138846
$ws.Range("D5").Value = 0.0958281402623046
$ws.Range("E5").Value = -0.6340838749161644
$ws.Range("F5").Value = 1.028095221473463
$ws.Range("G5").Value = 0.4088759330782192
$ws.Range("H5").Value = 0.6461796966979327
$ws.Range("I5").Value = 0.542902326822521
$ws.Range("J5").Value = -0.2486961005069136
$ws.Range("K5").Value = 0.3617928071605474
$ws.Range("B6").Value = 1.886257949266535
$ws.Range("C6").Value = 1.007855368014955
$ws.Range("D6").Value = 0.2779433528364856
$ws.Range("E6").Value = 1.940122449226113
$ws.Range("F6").Value = 1.320903160830869
$ws.Range("G6").Value = 1.558206924450583
$ws.Range("H6").Value = 1.454929554575171
$ws.Range("I6").Value = 0.6633311272457364
$ws.Range("J6").Value = 1.273820034913197
$ws.Range("K6").Value = 1.11229800409388
$ws.Range("B7").Value = 0.06388113204919779
$ws.Range("C7").Value = -0.6660308831292712
$ws.Range("D7").Value = 0.9961482132603562
$ws.Range("E7").Value = 0.3769289248651124
$ws.Range("F7").Value = 0.6142326884848259
$ws.Range("G7").Value = 0.5109553186094142
$ws.Range("H7").Value = -0.2806431087200204
$ws.Range("I7").Value = 0.3298457989474406
$ws.Range("J7").Value = 0.1683237681281231
$ws.Range("B8").Value = -0.6620629856161621
$ws.Range("C8").Value = 1.000116110773465
$ws.Range("D8").Value = 0.3808968223782215
$ws.Range("E8").Value = 0.6182005859979351
$ws.Range("F8").Value = 0.5149232161225235
$ws.Range("G8").Value = -0.2766752112069113
$ws.Range("H8").Value = 0.3338136964605497
$ws.Range("I8").Value = 0.1722916656412322
$ws.Range("B9").Value = 1.135756200943707
$ws.Range("C9").Value = 0.5165369125484629
$ws.Range("D9").Value = 0.7538406761681764
$ws.Range("E9").Value = 0.6505633062927647
$ws.Range("F9").Value = -0.1410351210366699
$ws.Range("G9").Value = 0.4694537866307911
$ws.Range("H9").Value = 0.3079317558114735
$ws.Range("B10").Value = 0.2757307306927982
$ws.Range("C10").Value = 0.5130344943125118
$ws.Range("D10").Value = 0.4097571244371001
$ws.Range("E10").Value = -0.3818413028923346
$ws.Range("F10").Value = 0.2286476047751264
$ws.Range("G10").Value = 0.06712557395580883
$ws.Range("B11").Value = 0.4677032790604154
$ws.Range("C11").Value = 0.3644259091850037
$ws.Range("D11").Value = -0.427172518144431
$ws.Range("E11").Value = 0.18331638952303
$ws.Range("F11").Value = 0.02179435870371246
$ws.Range("B12").Value = 0.2975644872489498
$ws.Range("C12").Value = -0.4940339400804848
$ws.Range("D12").Value = 0.1164549675869761
$ws.Range("E12").Value = -0.04506706323234141
$ws.Range("B13").Value = -0.5236201424372015
$ws.Range("C13").Value = 0.08686876523025952
$ws.Range("D13").Value = -0.07465326558905801
$ws.Range("B14").Value = 0.06871810850061863
$ws.Range("C14").Value = -0.0928039223186989
$ws.Range("B15").Value = -0.1108357465673982

# Clear the trailing cells that no longer have data for each row
$ws.Range("K7").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("B16").ClearContents()
